# Generate Report for Handoff
# Marks the ba931fbf-... and c3cf504a-... files as "Ready for handoff" on the
# Overview sheet and on the per-locale (zh-cn / de-de) detail sheets, refreshes
# the "latest handoff" timestamps, and records the "handback file not latest"
# error detail that the handoff-report generator produced for those two rows.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"

$ba931Error = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/794096d2b10aefbbbac36b31d4e6073cbb24a9aa/e2e/ba931fbf-3324-42eb-84fe-2fecb47adfef.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4d0aa4902bbf3af9d0f8dd3ab07c9fc523a1d810/e2e/ba931fbf-3324-42eb-84fe-2fecb47adfef.md."
$c3cfError  = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/794096d2b10aefbbbac36b31d4e6073cbb24a9aa/e2e/c3cf504a-c91b-4ff2-bfd7-6dea633c0f74.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4d0aa4902bbf3af9d0f8dd3ab07c9fc523a1d810/e2e/c3cf504a-c91b-4ff2-bfd7-6dea633c0f74.md."

# ---------------------------------------------------------------------------
# Sheet "Overview": rows 4 (ba931fbf...) and 5 (c3cf504a...)
#   columns: A File Name | B Path And Name | C Extension | D Publish URL
#            E zh-cn | F de-de | G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E4").Value = $statusReady
$overview.Range("F4").Value = $statusReady
$overview.Range("G4").Value = "2016-08-25 18:26:49"

$overview.Range("E5").Value = $statusReady
$overview.Range("F5").Value = $statusReady
$overview.Range("G5").Value = "2016-08-25 18:26:49"

# ---------------------------------------------------------------------------
# Sheet "zh-cn": rows 4 (ba931fbf...) and 5 (c3cf504a...)
#   columns: A Source File Name | B File Extension | C Status | D Source Path
#            E Priority | F Content Duplicate | G Latest Handoff File
#            H Latest Handoff Datetime | I Latest Target File
#            J Latest Handback File | K Latest Handback DateTime
#            L Reference Tokens | M To be localized | N Dependency From
#            O Has metadata | P Error Detail
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C4").Value = $statusReady
$zhcn.Range("H4").Value = "2016-08-25 18:26:45"
$zhcn.Range("P4").Value = $ba931Error

$zhcn.Range("C5").Value = $statusReady
$zhcn.Range("H5").Value = "2016-08-25 18:26:45"
$zhcn.Range("P5").Value = $c3cfError

# Error Detail column now holds long text -> widen it like the report does.
$zhcn.Columns.Item(16).ColumnWidth = 39.14

# ---------------------------------------------------------------------------
# Sheet "de-de": rows 4 (ba931fbf...) and 5 (c3cf504a...)
#   same column layout as "zh-cn"
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C4").Value = $statusReady
$dede.Range("H4").Value = "2016-08-25 18:26:49"
$dede.Range("P4").Value = $ba931Error

$dede.Range("C5").Value = $statusReady
$dede.Range("H5").Value = "2016-08-25 18:26:49"
$dede.Range("P5").Value = $c3cfError

$dede.Columns.Item(16).ColumnWidth = 39.14
